# Add new "Digit Destroyer" hint/feature strings to the language sheet.
# Cell values are written in the exact order needed so that the shared
# strings table grows with the same ordering as the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A95").Value = "digitDestroy"
$ws.Range("A96").Value = "digitDestroy_desc"
$ws.Range("B96").Value = "Select a blob to proceed."
$ws.Range("B97").Value = "Select a digit to destroy."
$ws.Range("A97").Value = "digitDestroy_modal_desc"
$ws.Range("A98").Value = "digitDestroy_dialog_1"
$ws.Range("A99").Value = "digitDestroy_dialog_2"
$ws.Range("B99").Value = "However, this will subtract from your score, so use it sparingly!"
$ws.Range("B95").Value = "DIGIT DESTROYER"
$ws.Range("B98").Value = "If you are having difficulty with certain numbers, press this button to remove some of its digits."

# Match the saved selection / scroll state from the edited workbook.
$ws.Range("B98").Select() | Out-Null
